$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the data rows (A2:D12) ascending by column A (time), keeping the
# header row (row 1) untouched.
$rng = $ws.Range("A2:D12")
$key = $ws.Range("A2:A12")
$rng.Sort($key, 1)
